$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append: date serial, first dose administered, second dose administered
$newRows = @(
    @(44536, 4517, 10652),
    @(44537, 4387, 10093),
    @(44538, 3786, 9339),
    @(44539, 3374, 9225),
    @(44540, 2859, 9250),
    @(44541, 2698, 8446),
    @(44542, 1193, 4268)
)

$startRow = 293
$lastExistingRow = $startRow - 1

$i = 0
foreach ($rowData in $newRows) {
    $r = $startRow + $i

    # Set the values first
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]

    $i = $i + 1
}

$lastRow = $startRow + $newRows.Count - 1

# Copy the date-column formatting from the previous final row down onto the
# newly added date cells so they keep the existing dd/mm/yyyy number format.
$srcRange = "A" + $lastExistingRow
$destRange = "A" + $startRow + ":A" + $lastRow
$ws.Range($srcRange).Copy()
$ws.Range($destRange).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the saved selection to match the workbook's new selection
$ws.Range("C298").Select()
